$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-09-23 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-09-24 Wednesday", 2) | Out-Null
$d.Content.Find.Execute("869÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "651÷7=", 2) | Out-Null
$d.Content.Find.Execute("792÷2=", $true, $false, $false, $false, $false, $true, 1, $false, "775÷9=", 2) | Out-Null
$d.Content.Find.Execute("892÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "418÷9=", 2) | Out-Null
$d.Content.Find.Execute("566÷3=", $true, $false, $false, $false, $false, $true, 1, $false, "632÷6=", 2) | Out-Null
$d.Content.Find.Execute("195÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "185÷7=", 2) | Out-Null
$d.Content.Find.Execute("470÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "376÷8=", 2) | Out-Null
$d.Content.Find.Execute("298÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "883÷9=", 2) | Out-Null
$d.Content.Find.Execute("478÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "863÷3=", 2) | Out-Null
$d.Content.Find.Execute("389÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "570÷2=", 2) | Out-Null
$d.Content.Find.Execute("759÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "272÷2=", 2) | Out-Null
$d.Content.Find.Execute("591÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "263÷8=", 2) | Out-Null
$d.Content.Find.Execute("841÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "569÷7=", 2) | Out-Null
$d.Content.Find.Execute("387÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "639÷5=", 2) | Out-Null
$d.Content.Find.Execute("607÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "280÷2=", 2) | Out-Null
$d.Content.Find.Execute("390÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "623÷2=", 2) | Out-Null
$d.Content.Find.Execute("210÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "625÷4=", 2) | Out-Null
$d.Content.Find.Execute("633÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "552÷3=", 2) | Out-Null
$d.Content.Find.Execute("343÷8=", $true, $false, $false, $false, $false, $true, 1, $false, "541÷8=", 2) | Out-Null
$d.Content.Find.Execute("813÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "272÷9=", 2) | Out-Null
$d.Content.Find.Execute("164÷5=", $true, $false, $false, $false, $false, $true, 1, $false, "885÷9=", 2) | Out-Null
$d.Content.Find.Execute("811÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "757÷3=", 2) | Out-Null
$d.Content.Find.Execute("139÷7=", $true, $false, $false, $false, $false, $true, 1, $false, "368÷6=", 2) | Out-Null
$d.Content.Find.Execute("519÷6=", $true, $false, $false, $false, $false, $true, 1, $false, "742÷4=", 2) | Out-Null
$d.Content.Find.Execute("660÷4=", $true, $false, $false, $false, $false, $true, 1, $false, "564÷3=", 2) | Out-Null
$d.Content.Find.Execute("814÷9=", $true, $false, $false, $false, $false, $true, 1, $false, "705÷7=", 2) | Out-Null
